$wb = $excel.ActiveWorkbook

# --- Measures sheet ---
$wsM = $wb.Worksheets.Item("Measures")

$wsM.Range("B2").Value = "`nSUM(`n    PNP_DELIVERIES_SSC_001[PnP Returns Qty CU]`n)"
$wsM.Range("E2").Value = "This calculation would add the number of return products in a given unit, such as boxes or pallets, for the selected PnP Returns Qty field and provide a total for the selected period."
$wsM.Range("B3").Value = "`nSUM( PNP_DELIVERIES_SSC_001[PnP Qty Received CU] )"
$wsM.Range("E3").Value = "This calculation sums up the total quantity of product received by a company in the `"PnP Qty Received CU`" column of the table `"PNP_DELIVERIES_SSC_001`". This calculation is useful for understanding how much of a product a business has received from outside suppliers or distribution"
$wsM.Range("B4").Value = "`nCALCULATE(`n    DIVIDE(`n        PNP_DELIVERIES_SSC_001[Quantity Returned], `n        PNP_DELIVERIES_SSC_001[Quantity Recieved], `n        BLANK()`n    ) `n)*100"
$wsM.Range("E4").Value = "This calculation divides the quantity of items returned by the quantity of items received, and then multiplies the result by 100. This calculation can help determine the return rate as a percentage of all items received."
$wsM.Range("B5").Value = "`n-LASTDATE(Ax[Date])`nTODAY() - LASTDATE(Ax[Date])"
$wsM.Range("E5").Value = "- 7`nThis calculation subtracts 7 days from the current date using the TODAY() function. This can be used to determine a date 7 days in the past from the current date. This is useful if you need to compare and analyze data within a one week period."
$wsM.Range("B6").Value = "`nCALCULATE(`n    PNP_DELIVERIES_SSC_001[Quantity Returned],`n    PNP_DELIVERIES_SSC_001[Prod Sub Category] = `"Premium`"`n)"
$wsM.Range("E6").Value = "This calculation is finding the total quantity of `"Premium`" products that were returned. It is looking through the data set PNP_DELIVERIES_SSC_001 and searching for the sum of the quantity returned for every product that falls under the `"Premium`" product subcategory."
$wsM.Range("B7").Value = "`nCALCULATE(`n    PNP_DELIVERIES_SSC_001[Quantity Returned], `n    PNP_DELIVERIES_SSC_001[Prod Sub Category] = `"Standard`"`n)"
$wsM.Range("E7").Value = "This calculation is finding the sum of the `"Quantity Returned`" from the table `"PNP_DELIVERIES_SSC_001`" where the `"Prod Sub Category`" is equal to `"Standard`"."
$wsM.Range("E8").Value = "This calculation is for the total quantity returned for all products in the `"Dumpy`" subcategory in the PNP_DELIVERIES_SSC_001 table. It is used to determine the total quantity of defective products or products returned for a specific product subcategory."
$wsM.Range("B9").Value = "`nCALCULATE(`n    PNP_DELIVERIES_SSC_001[Quantity Returned], `n    PNP_DELIVERIES_SSC_001[Prod Sub Category] = `"Smart`"`n)"
$wsM.Range("E9").Value = "This calculation finds the total quantity of products returned in the `"Smart`" product subcategory of the PNP_DELIVERIES_SSC_001 data set. It multiplies each 'Quantity Returned' value in the data set by a filter on 'Prod Sub Category'. This produces the total"
$wsM.Range("B10").Value = "`nCALCULATE(`n    PNP_DELIVERIES_SSC_001[Quantity Returned],`n    PNP_DELIVERIES_SSC_001[Prod Sub Category] = `"Buns & Rolls`"`n)"
$wsM.Range("E10").Value = "This calculation is finding the total quantity returned for products in the sub category `"Buns & Rolls`" in the PNP_DELIVERIES_SSC_001 table."
$wsM.Range("B11").Value = "`nCALCULATE(`n    PNP_DELIVERIES_SSC_001[Quantity Returned],`n    PNP_DELIVERIES_SSC_001[Prod Sub Category] = `"Everyday+`"`n)"
$wsM.Range("E11").Value = "This calculation determines the total quantity of items returned for the product subcategory `"Everyday+`". It looks through the entries in the table PNP_DELIVERIES_SSC_001 to count up how many items in the Everyday+ subcategory have been returned."
$wsM.Range("E12").Value = "This calculation returns the greatest date among the 'Date' column of the 'PNP_DELIVERIES_SSC_001' table. It will be the most recent date among all 'Date' values in the column. This could be used to identify the time of the latest delivery in the table"

# --- Source Information sheet ---
$wsS = $wb.Worksheets.Item("Source Information")
$wsS.Range("A2").Value = 1
$wsS.Range("I2").Value = "1. This renames the `"StoreType`" column to `"Store Type`" in the source table.`n`n"

# Resize the Source table to match the actual data range (A1:I2)
$lo = $wsS.ListObjects.Item(1)
$lo.Resize($wsS.Range("A1:I2"))

Write-Host "done"
